$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2,6).Value2 = 204
$ws1.Cells.Item(5,6).Value2 = 179
$ws1.Cells.Item(6,6).Value2 = 536
$ws1.Cells.Item(7,6).Value2 = 48
$ws1.Cells.Item(8,6).Value2 = 9665
$ws1.Cells.Item(9,6).Value2 = 50
$ws1.Cells.Item(10,6).Value2 = 2603
$ws1.Cells.Item(11,6).Value2 = 203
$ws1.Cells.Item(12,6).Value2 = 2378
$ws1.Cells.Item(13,6).Value2 = 2618
$ws1.Cells.Item(14,6).Value2 = 1398
$ws1.Cells.Item(16,6).Value2 = 2046
$ws1.Cells.Item(17,6).Value2 = 43
$ws1.Cells.Item(18,6).Value2 = 73
$ws1.Cells.Item(19,6).Value2 = 360
$ws1.Cells.Item(21,6).Value2 = 61
$ws1.Cells.Item(22,6).Value2 = 293
$ws1.Cells.Item(23,6).Value2 = 57
$ws1.Cells.Item(24,6).Value2 = 129
$ws1.Cells.Item(26,6).Value2 = 1268
$ws1.Cells.Item(29,6).Value2 = 116
$ws1.Cells.Item(31,6).Value2 = 1639
$ws1.Cells.Item(32,6).Value2 = 2735
$ws1.Cells.Item(34,6).Value2 = 970
$ws1.Cells.Item(35,6).Value2 = 340
$ws1.Cells.Item(37,6).Value2 = 34
$ws1.Cells.Item(38,6).Value2 = 45

# ---- Sheet: 演出 (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")

# Rows 4-14 take on the (shifted) content of what used to be rows 6-16;
# column A (the running index) is left untouched.
# -- row 4 --
$c = $ws2.Cells.Item(4,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-08-19"
$ws2.Cells.Item(4,3).Value2 = "北京·音阅派国漫演唱会-《一人之下》动画八周年专场演唱会"
$ws2.Cells.Item(4,4).Value2 = "中关村南大街33号中国国家图书馆内 国图艺术中心"
$ws2.Cells.Item(4,5).Value2 = "2024.08.19 19:30-08.19 21:00"
$ws2.Cells.Item(4,6).Value2 = 161
$ws2.Cells.Item(4,7).Value2 = 380
$ws2.Cells.Item(4,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89052"
$ws2.Cells.Item(4,9).Value2 = "//i2.hdslb.com/bfs/openplatform/202407/wtYvGYL51720603864335.png"
# -- row 5 --
$c = $ws2.Cells.Item(5,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-08-23"
$ws2.Cells.Item(5,3).Value2 = "北京·《山丘》音乐教父 经典情歌金曲翻唱演唱会"
$ws2.Cells.Item(5,4).Value2 = "大江胡同121号2幢负1层 北京门空间 TheDoorLiveHouse"
$ws2.Cells.Item(5,5).Value2 = "2024.08.23 19:30-08.23 21:00"
$ws2.Cells.Item(5,6).Value2 = 1
$ws2.Cells.Item(5,7).Value2 = 98
$ws2.Cells.Item(5,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89358"
$ws2.Cells.Item(5,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202407/noqwx8Qu1721116074567.jpeg"
# -- row 6 --
$c = $ws2.Cells.Item(6,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-08-24"
$ws2.Cells.Item(6,3).Value2 = "北京·最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会"
$ws2.Cells.Item(6,4).Value2 = "亮马桥路40号(近好运街) 北京世纪剧院"
$ws2.Cells.Item(6,5).Value2 = "2024.08.24 19:30-08.24 21:00"
$ws2.Cells.Item(6,6).Value2 = 14
$ws2.Cells.Item(6,7).Value2 = 153
$ws2.Cells.Item(6,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=86217"
$ws2.Cells.Item(6,9).Value2 = "//i2.hdslb.com/bfs/openplatform/202405/BDyblKrJ1716427731729.jpeg"
# -- row 7 --
$c = $ws2.Cells.Item(7,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-09-30"
$ws2.Cells.Item(7,3).Value2 = "北京·VGL 中国巡演 2024 VIDEO GAME LIVE 魔兽世界音乐会"
$ws2.Cells.Item(7,4).Value2 = "西直门外大街135号  北展剧场"
$ws2.Cells.Item(7,5).Value2 = "2024.09.30 19:30-09.30 21:30"
$ws2.Cells.Item(7,6).Value2 = 8
$ws2.Cells.Item(7,7).Value2 = 180
$ws2.Cells.Item(7,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89144"
$ws2.Cells.Item(7,9).Value2 = "//i2.hdslb.com/bfs/openplatform/202407/5YIwe8lU1720605586333.jpeg"
# -- row 8 --
$c = $ws2.Cells.Item(8,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-10-01"
$ws2.Cells.Item(8,3).Value2 = "北京·VGL 中国巡演 2024 VIDEO GAMES LIVE 暴雪游戏音乐会"
$ws2.Cells.Item(8,4).Value2 = "西直门外大街135号  北展剧场"
$ws2.Cells.Item(8,5).Value2 = "2024.10.01 19:30-10.01 21:30"
$ws2.Cells.Item(8,6).Value2 = 15
$ws2.Cells.Item(8,7).Value2 = 180
$ws2.Cells.Item(8,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89083"
$ws2.Cells.Item(8,9).Value2 = "//i1.hdslb.com/bfs/openplatform/202407/yMoDGuXs1720607500874.jpeg"
# -- row 9 --
$c = $ws2.Cells.Item(9,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-10-10"
$ws2.Cells.Item(9,3).Value2 = "北京·黑白键上的音乐地图——孩子们的钢琴协奏曲之夜"
$ws2.Cells.Item(9,4).Value2 = "北新华街1号 北京音乐厅"
$ws2.Cells.Item(9,5).Value2 = "2024.10.10 19:30-10.10 21:00"
$ws2.Cells.Item(9,6).Value2 = 1
$ws2.Cells.Item(9,7).Value2 = 153
$ws2.Cells.Item(9,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=86881"
$ws2.Cells.Item(9,9).Value2 = "//i1.hdslb.com/bfs/openplatform/202406/K3oihoH91717474488019.jpeg"
# -- row 10 --
$c = $ws2.Cells.Item(10,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-10-11"
$ws2.Cells.Item(10,3).Value2 = "北京·官方唯一授权-周杰伦作品国风音乐会"
$ws2.Cells.Item(10,4).Value2 = "西直门外大街135号  北展剧场"
$ws2.Cells.Item(10,5).Value2 = "2024.10.11 19:30-10.11 21:00"
$ws2.Cells.Item(10,6).Value2 = 11
$ws2.Cells.Item(10,7).Value2 = 126
$ws2.Cells.Item(10,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=88666"
$ws2.Cells.Item(10,9).Value2 = "//i1.hdslb.com/bfs/openplatform/202407/2KgWinEn1720077808243.jpeg"
# -- row 11 --
$c = $ws2.Cells.Item(11,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-10-25"
$ws2.Cells.Item(11,3).Value2 = "北京·伦敦西区音乐剧明星演唱会-经典版"
$ws2.Cells.Item(11,4).Value2 = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$ws2.Cells.Item(11,5).Value2 = "2024.10.25 19:30-10.26 21:30"
$ws2.Cells.Item(11,6).Value2 = 4
$ws2.Cells.Item(11,7).Value2 = 144
$ws2.Cells.Item(11,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89359"
$ws2.Cells.Item(11,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202407/PzPiEKUI1721114840552.jpeg"
# -- row 12 --
$c = $ws2.Cells.Item(12,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-10-26"
$ws2.Cells.Item(12,3).Value2 = "北京·伦敦西区音乐剧明星演唱会（摇滚版）"
$ws2.Cells.Item(12,4).Value2 = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$ws2.Cells.Item(12,5).Value2 = "2024.10.26 14:30-10.26 16:30"
$ws2.Cells.Item(12,6).Value2 = 5
$ws2.Cells.Item(12,7).Value2 = 144
$ws2.Cells.Item(12,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89400"
$ws2.Cells.Item(12,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202407/TYPRpfu21721116217467.jpeg"
# -- row 13 --
$c = $ws2.Cells.Item(13,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-10-26"
$ws2.Cells.Item(13,3).Value2 = "北京·变形金刚音乐会40周年特别版"
$ws2.Cells.Item(13,4).Value2 = "中关村南大街33号国家图书馆北门 国图艺术中心音乐厅"
$ws2.Cells.Item(13,5).Value2 = "2024.10.26 19:30-10.26 21:30"
$ws2.Cells.Item(13,6).Value2 = 33
$ws2.Cells.Item(13,7).Value2 = 171
$ws2.Cells.Item(13,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89970"
$ws2.Cells.Item(13,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202407/TwvRQI041722150343639.jpeg"
# -- row 14 --
$c = $ws2.Cells.Item(14,2)
$c.NumberFormat = "@"
$c.Value2 = "2024-11-30"
$ws2.Cells.Item(14,3).Value2 = "北京·花たん 2024 LIVE in Beijing"
$ws2.Cells.Item(14,4).Value2 = "复兴路69号院2号136、G23室 Mao Livehouse北京五棵松店"
$ws2.Cells.Item(14,5).Value2 = "2024.11.30 14:00-11.30 15:30"
$ws2.Cells.Item(14,6).Value2 = 147
$ws2.Cells.Item(14,7).Value2 = 380
$ws2.Cells.Item(14,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=90341"
$ws2.Cells.Item(14,9).Value2 = "//i2.hdslb.com/bfs/openplatform/202408/wfGEn3sY1722910561352.jpeg"

# Remove the now-obsolete rows 15 and 16 (events dropped from the feed)
$ws2.Rows.Item(15).Resize(2).Delete() | Out-Null

# ---- Sheet: 本地生活 (Local Life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2,6).Value2 = 719
$ws3.Cells.Item(3,6).Value2 = 940

# ---- Sheet: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2,6).Value2 = 204
$ws4.Cells.Item(3,6).Value2 = 719
$ws4.Cells.Item(4,6).Value2 = 940
$ws4.Cells.Item(9,6).Value2 = 179
$ws4.Cells.Item(10,6).Value2 = 536
$ws4.Cells.Item(11,6).Value2 = 48
$ws4.Cells.Item(12,6).Value2 = 9665
$ws4.Cells.Item(14,6).Value2 = 50
$ws4.Cells.Item(15,6).Value2 = 2603
$ws4.Cells.Item(16,6).Value2 = 203
$ws4.Cells.Item(17,6).Value2 = 2378
$ws4.Cells.Item(18,6).Value2 = 2618
$ws4.Cells.Item(21,6).Value2 = 2046
$ws4.Cells.Item(22,6).Value2 = 43
$ws4.Cells.Item(23,6).Value2 = 73
$ws4.Cells.Item(24,6).Value2 = 360
$ws4.Cells.Item(26,6).Value2 = 61
$ws4.Cells.Item(27,6).Value2 = 293
$ws4.Cells.Item(28,6).Value2 = 57
$ws4.Cells.Item(29,6).Value2 = 129
$ws4.Cells.Item(31,6).Value2 = 1268
$ws4.Cells.Item(34,6).Value2 = 116
$ws4.Cells.Item(36,6).Value2 = 1639
$ws4.Cells.Item(38,6).Value2 = 2735
$ws4.Cells.Item(39,6).Value2 = 970
$ws4.Cells.Item(41,6).Value2 = 340
$ws4.Cells.Item(46,6).Value2 = 45
$ws4.Cells.Item(49,6).Value2 = 147
$ws4.Cells.Item(50,6).Value2 = 147
